# إضافة حدث جديد في Card17
# Adds a new service-log row (row 17) to the "Card17" sheet, and fills in
# the previously blank B:K cells of the prior last row (row 16) with the
# literal "nan" placeholder text used throughout the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card17")

# Row 16 (previous last row): its measurement columns (B:K) were blank;
# they now carry the sheet's standard "nan" placeholder text.
foreach ($col in @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")) {
    $ws.Range($col + "16").Value = "nan"
}

# New row 17: a new service/event entry for this card.
# A17 mirrors the "card" number convention used by every other row in the
# sheet, which is stored as text rather than a number.
$ws.Range("A17").Value = "'17"
$ws.Range("L17").Value = "14\8\2025"
$ws.Range("N17").Value = "تم تغيير زيت الجيربوكس"
$ws.Range("O17").Value = "تيم العمل"
